# Auto-generated: apply cell-value updates per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2092.5264
$ws.Range("I51").Value = 2015.4546
$ws.Range("J51").Value = 2198.5
$ws.Range("K51").Value = 2015.4546
$ws.Range("L51").Value = 2198.5
$ws.Range("M51").Value = -1531.4546
$ws.Range("N51").Value = -3166.5
$ws.Range("H121").Value = 749.8570999999999
$ws.Range("J121").Value = 746
$ws.Range("L121").Value = 2238
$ws.Range("N121").Value = -5732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 34224.453
$ws.Range("I2").Value = 38394.742
$ws.Range("J2").Value = 6075
$ws.Range("K2").Value = 38394.742
$ws.Range("L2").Value = 6075
$ws.Range("M2").Value = -38281.742
$ws.Range("N2").Value = -6301
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -84
$ws.Range("H32").Value = 19235.441
$ws.Range("I32").Value = 16367.864
$ws.Range("K32").Value = 16367.864
$ws.Range("M32").Value = -16080.864
$ws.Range("H63").Value = 2105.5557
$ws.Range("I63").Value = 2105.5557
$ws.Range("K63").Value = 2105.5557
$ws.Range("M63").Value = -1419.5557
$ws.Range("H66").Value = 2105.5557
$ws.Range("I66").Value = 2105.5557
$ws.Range("K66").Value = 10527.7785
$ws.Range("M66").Value = -7095.7785
$ws.Range("H76").Value = 50288
$ws.Range("J76").Value = 50288
$ws.Range("L76").Value = 50288
$ws.Range("N76").Value = -50964
$ws.Range("H79").Value = 50288
$ws.Range("J79").Value = 50288
$ws.Range("L79").Value = 50288
$ws.Range("N79").Value = -52628
$ws.Range("H116").Value = 34224.453
$ws.Range("I116").Value = 38394.742
$ws.Range("J116").Value = 6075
$ws.Range("K116").Value = 38394.742
$ws.Range("L116").Value = 6075
$ws.Range("M116").Value = -36100.742
$ws.Range("N116").Value = -10663

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 34224.453
$ws.Range("I3").Value = 38394.742
$ws.Range("J3").Value = 6075
$ws.Range("K3").Value = 38394.742
$ws.Range("L3").Value = 6075
$ws.Range("M3").Value = -38280.742
$ws.Range("N3").Value = -6303
$ws.Range("H20").Value = 47651030
$ws.Range("I20").Value = 162626.5
$ws.Range("J20").Value = 58824776
$ws.Range("K20").Value = 162626.5
$ws.Range("L20").Value = 58824776
$ws.Range("M20").Value = -162379.5
$ws.Range("N20").Value = -58825270
$ws.Range("H54").Value = 1441.4445
$ws.Range("I54").Value = 1567.5714
$ws.Range("J54").Value = 1000
$ws.Range("K54").Value = 1567.5714
$ws.Range("L54").Value = 1000
$ws.Range("M54").Value = -1083.5714
$ws.Range("N54").Value = -1968
$ws.Range("H99").Value = 62502096
$ws.Range("I99").Value = 83335210
$ws.Range("J99").Value = 2755.5
$ws.Range("K99").Value = 83335210
$ws.Range("L99").Value = 2755.5
$ws.Range("M99").Value = -83333712
$ws.Range("N99").Value = -5751.5
$ws.Range("H134").Value = 27442.41
$ws.Range("I134").Value = 1906.16
$ws.Range("J134").Value = 73042.86
$ws.Range("K134").Value = 5718.48
$ws.Range("L134").Value = 219128.58
$ws.Range("M134").Value = -3183.48
$ws.Range("N134").Value = -224198.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2400.0435
$ws.Range("I31").Value = 944.55554
$ws.Range("J31").Value = 3335.7144
$ws.Range("K31").Value = 944.55554
$ws.Range("L31").Value = 3335.7144
$ws.Range("M31").Value = -649.55554
$ws.Range("N31").Value = -3925.7144
$ws.Range("H34").Value = 2400.0435
$ws.Range("I34").Value = 944.55554
$ws.Range("J34").Value = 3335.7144
$ws.Range("K34").Value = 944.55554
$ws.Range("L34").Value = 3335.7144
$ws.Range("M34").Value = -742.55554
$ws.Range("N34").Value = -3739.7144
$ws.Range("H50").Value = 16165.333
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 16165.333
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 16165.333
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -17415.333
$ws.Range("H51").Value = 18725
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18725
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18725
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -20197
$ws.Range("H58").Value = 5256.2964
$ws.Range("I58").Value = 1034.6666
$ws.Range("J58").Value = 7367.1113
$ws.Range("K58").Value = 1034.6666
$ws.Range("L58").Value = 7367.1113
$ws.Range("M58").Value = -831.6666
$ws.Range("N58").Value = -7773.1113
$ws.Range("H61").Value = 18725
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 18725
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 18725
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -19421
$ws.Range("H136").Value = 5256.2964
$ws.Range("I136").Value = 1034.6666
$ws.Range("J136").Value = 7367.1113
$ws.Range("K136").Value = 3103.9998
$ws.Range("L136").Value = 22101.3339
$ws.Range("M136").Value = -553.9998000000001
$ws.Range("N136").Value = -27201.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 6912
$ws.Range("I118").Value = 1500
$ws.Range("J118").Value = 7453.2
$ws.Range("K118").Value = 4500
$ws.Range("L118").Value = 22359.6
$ws.Range("M118").Value = -3257
$ws.Range("N118").Value = -24845.6
$ws.Range("H122").Value = 440.76923
$ws.Range("I122").Value = 383.38095
$ws.Range("J122").Value = 681.8
$ws.Range("K122").Value = 3450.42855
$ws.Range("L122").Value = 6136.2
$ws.Range("M122").Value = -1000.42855
$ws.Range("N122").Value = -11036.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2745.5454
$ws.Range("I97").Value = 2819
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 2819
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -2323
$ws.Range("N97").Value = -3003
$ws.Range("H132").Value = 1982.5
$ws.Range("I132").Value = 1399.0555
$ws.Range("J132").Value = 3732.8333
$ws.Range("K132").Value = 4197.166499999999
$ws.Range("L132").Value = 11198.4999
$ws.Range("M132").Value = -1667.166499999999
$ws.Range("N132").Value = -16258.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 8710.5
$ws.Range("J42").Value = 8710.5
$ws.Range("L42").Value = 8710.5
$ws.Range("N42").Value = -9836.5
$ws.Range("H46").Value = 451.92307
$ws.Range("J46").Value = 483.57144
$ws.Range("L46").Value = 483.57144
$ws.Range("N46").Value = -859.5714399999999
$ws.Range("H49").Value = 8710.5
$ws.Range("J49").Value = 8710.5
$ws.Range("L49").Value = 8710.5
$ws.Range("N49").Value = -9004.5
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232
$ws.Range("H100").Value = 5907891
$ws.Range("I100").Value = 7015301.5
$ws.Range("J100").Value = 1700.6666
$ws.Range("K100").Value = 7015301.5
$ws.Range("L100").Value = 1700.6666
$ws.Range("M100").Value = -7014760.5
$ws.Range("N100").Value = -2782.6666
$ws.Range("H122").Value = 7908.722
$ws.Range("I122").Value = 10417.667
$ws.Range("J122").Value = 2890.8333
$ws.Range("K122").Value = 31253.001
$ws.Range("L122").Value = 8672.499899999999
$ws.Range("M122").Value = -28803.001
$ws.Range("N122").Value = -13572.4999
$ws.Range("H132").Value = 1574657.6
$ws.Range("I132").Value = 1986207.5
$ws.Range("J132").Value = 3285.182
$ws.Range("K132").Value = 5958622.5
$ws.Range("L132").Value = 9855.545999999998
$ws.Range("M132").Value = -5956092.5
$ws.Range("N132").Value = -14915.546
$ws.Range("H136").Value = 1785.3606
$ws.Range("I136").Value = 1025.2667
$ws.Range("J136").Value = 2520.9355
$ws.Range("K136").Value = 3075.800099999999
$ws.Range("L136").Value = 7562.806500000001
$ws.Range("M136").Value = -525.8000999999995
$ws.Range("N136").Value = -12662.8065

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 604.8333
$ws.Range("I113").Value = 709.6667
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 2129.0001
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 40.9998999999998
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 1196.2727
$ws.Range("I122").Value = 1054
$ws.Range("J122").Value = 1680
$ws.Range("K122").Value = 3162
$ws.Range("L122").Value = 5040
$ws.Range("M122").Value = -712
$ws.Range("N122").Value = -9940
